$wb = $excel.ActiveWorkbook

# Add the new "negative_score" worksheet after the last existing sheet (designation)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "negative_score"

# Header row
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "numeric"

# Weightage row
$ws.Range("A2").Value = "weightage"
$ws.Range("B2").Value = 1

# Score rows
$ws.Range("A3").Value = "competitor"
$ws.Range("B3").Value = -100

$ws.Range("A4").Value = "no competitor"
$ws.Range("B4").Value = 0

# Selection on the new sheet matches B5 (one below the data, like the original workbook)
$ws.Range("B5").Select() | Out-Null
